# Update marksheet figures: correct/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Marking" row - right answers count
$ws.Range("B11").Value = 5

# "Total" row - total marks obtained and the "obtained/max" text
$ws.Range("B12").Value = 110
$ws.Range("E12").Value = "110/140"

$wb.Save()
